# Inserts a new weekly price record at row 213 of the "Apio" sheet.
# Every existing row from 213 down to 252 shifts down by one row
# (to 214..253); the newly opened row 213 receives the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 213:252 down to 214:253, opening up a blank row 213.
$ws.Rows("213:213").Insert()

# Populate the newly inserted row 213 with the new record.
$ws.Range("A213").Value = 4
$ws.Range("B213").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C213").Value = "Los Lagos"
$ws.Range("D213").Value = 44711
$ws.Range("E213").Value = 10
$ws.Range("F213").Value = 100112017
$ws.Range("G213").Value = "Apio"
$ws.Range("H213").Value = "Americana (o)"
$ws.Range("I213").Value = "Primera"
$ws.Range("J213").Value = 25
$ws.Range("K213").Value = 12000
$ws.Range("L213").Value = 12000
$ws.Range("M213").Value = 12000
$ws.Range("N213").Value = "$/docena de matas"
$ws.Range("O213").Value = "Región de Coquimbo"
$ws.Range("P213").Value = 2000
$ws.Range("Q213").Value = 6
$ws.Range("R213").Value = "Hortaliza"
